$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - mirrors style/formatting of existing header cells (e.g. G1)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New "Save" column values
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
